$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("diffusion -> protein mpnn -> AF"): insert 3 new columns (E:G) for AF3 predictions ---
$ws1.Columns("E:G").Insert()

# Headers for the new columns
$ws1.Range("E1").Value = "overall confidence"
$ws1.Range("F1").Value = "ligand score"
$ws1.Range("G1").Value = "seq rec"

# Match the formatting (centered style) already used by column D
$ws1.Range("E1:G4").HorizontalAlignment = -4108
$ws1.Range("E1:G4").ColumnWidth = $ws1.Range("D1").ColumnWidth

# New AF3 prediction values
$ws1.Range("E2").Value = 0.40350000000000003
$ws1.Range("F2").Value = 0.40350000000000003
$ws1.Range("G2").Value = 0.31940000000000002

$ws1.Range("E3").Value = 0.43209999999999998
$ws1.Range("F3").Value = 0.43209999999999998
$ws1.Range("G3").Value = 0.36630000000000001

$ws1.Range("E4").Value = 0.45250000000000001
$ws1.Range("F4").Value = 0.45250000000000001
$ws1.Range("G4").Value = 0.45950000000000002

# --- Sheet2 ("ligand mpnn -> AF3 -> docking"): unify the group-label cell
#     formatting for A8:A10 with the A2:A4 / A5:A7 groups (drop the stray border) ---
$ws2.Range("A2").Copy()
$ws2.Range("A8:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet2 becomes the active sheet/tab ---
$ws2.Activate()
$ws2.Range("B13").Select()

# Re-select sheet1's new selection anchor before giving focus back to sheet2
$ws1.Range("G5").Select()
$ws2.Activate()
